# T1432 - Activity+new jon TAS Job Type
# Adds a new FAS job-type mapping row to the "AddOpportunity" sheet and a
# matching "Engagement" sheet row for the new "TAS" job type conversion.

$wb = $excel.ActiveWorkbook

$wsOpportunity = $wb.Worksheets.Item("AddOpportunity")
$wsEngagement  = $wb.Worksheets.Item("Engagement")
$wsContact     = $wb.Worksheets.Item("AddContact")

# --- AddOpportunity: duplicate row 12 into row 13, changing only the
# FASJobType (column C) value to the new job type.
$wsOpportunity.Range("A12:AB12").Copy()
$wsOpportunity.Range("A13:AB13").PasteSpecial(-4122) # xlPasteFormats
$wsOpportunity.Range("A13:AB13").Value2 = $wsOpportunity.Range("A12:AB12").Value2
$wsOpportunity.Range("C13").Value2 = "TAS - ESG Due Diligence & Analytics"
$wsOpportunity.Range("AD13").Value2 = $wsOpportunity.Range("AD12").Value2
$wsOpportunity.Rows.Item(13).RowHeight = $wsOpportunity.Rows.Item(12).RowHeight

# --- Engagement: add the matching conversion row (copy formats from the
# existing rows that carry the right style for each column).
$wsEngagement.Range("A12").Copy()
$wsEngagement.Range("A13").PasteSpecial(-4122) # xlPasteFormats
$wsEngagement.Range("A13").Value2 = "Retained"

$wsEngagement.Range("B2").Copy()
$wsEngagement.Range("B13").PasteSpecial(-4122) # xlPasteFormats
$wsEngagement.Range("B13").Value2 = "TAS - Due Diligence Services"

$wsEngagement.Range("C12").Copy()
$wsEngagement.Range("C13").PasteSpecial(-4122) # xlPasteFormats
$wsEngagement.Range("C13").Value2 = "HL Capital, Inc."

# --- Restore / update selections to match the saved view state.
$wsOpportunity.Activate()
$wsOpportunity.Rows.Item(2).Select()

$wsContact.Activate()
$wsContact.Range("F25").Select()

$wsEngagement.Activate()
$wsEngagement.Range("D12").Select()
